$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append " (Chrome)" to the existing bug descriptions in rows 13-16 ---
$ws.Range("A13").Value = $ws.Range("A13").Value2 + " (Chrome)"
$ws.Range("A14").Value = $ws.Range("A14").Value2 + " (Chrome)"
$ws.Range("A15").Value = $ws.Range("A15").Value2 + " (Chrome)"
$ws.Range("A16").Value = $ws.Range("A16").Value2 + " (Chrome)"

# Row 16 grows by one wrapped line once " (Chrome)" is appended
$ws.Rows.Item(16).RowHeight = 187.2

# --- Add three new rows of testing data (17-19), using row 12 as the
#     formatting template (same column layout: A/B/E plain, C "priority"
#     style, D date style) ---
$ws.Range("A12:I12").Copy()
$ws.Range("A17:I19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A17").Value = "if you add a section and then click on the textbox which says section title, if you press enter, every added section disappears (Chrome)"
$ws.Range("B17").Value = "edit.html"
$ws.Range("C17").Value = "Medium"
$ws.Range("D17").Value = 43089
$ws.Range("E17").Value = "Xavier Kuttamparambil"
$ws.Rows.Item(17).RowHeight = 93.6

$ws.Range("A18").Value = "same as above for the textbox diagram url (Chrome)"
$ws.Range("B18").Value = "edit.html"
$ws.Range("C18").Value = "Medium"
$ws.Range("D18").Value = 43089
$ws.Range("E18").Value = "Xavier Kuttamparambil"
$ws.Rows.Item(18).RowHeight = 46.8

$ws.Range("A19").Value = "the add image button also makes every section disappear as well (Chrome)"
$ws.Range("B19").Value = "edit.html"
$ws.Range("C19").Value = "Medium"
$ws.Range("D19").Value = 43090
$ws.Range("E19").Value = "Xavier Kuttamparambil"
$ws.Rows.Item(19).RowHeight = 62.4

# --- Restore the view state to what it was left at after the edits ---
$ws.Range("D18").Select()
